$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing 2014 count value
$ws.Range("B5").Value = 432

# Add new row for 2015 (force A6 to be stored as text "2015", not a number)
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2015"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = 292
